# Auto-generated edit script applying odds updates per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G11").Value = 3.4
$ws.Range("I11").Value = 2.15
$ws.Range("J11").Value = 4
$ws.Range("L11").Value = 2.77
$ws.Range("P11").Value = 2.45
$ws.Range("Q11").Value = 2.22
$ws.Range("S11").Value = 1.47
$ws.Range("T11").Value = 2.35
$ws.Range("W11").Value = 8
$ws.Range("X11").Value = 16.5
$ws.Range("Y11").Value = 12.5
$ws.Range("Z11").Value = 50
$ws.Range("AA11").Value = 37
$ws.Range("AB11").Value = 50
$ws.Range("AG11").Value = 6.1
$ws.Range("AH11").Value = 9.25
$ws.Range("AI11").Value = 9.25
$ws.Range("AJ11").Value = 20
$ws.Range("AK11").Value = 20
$ws.Range("AL11").Value = 37
$ws.Range("AN11").Value = 5.1
$ws.Range("AO11").Value = 19.5
$ws.Range("AP11").Value = 29
$ws.Range("AQ11").Value = 110
$ws.Range("AS11").Value = 450
$ws.Range("AW11").Value = 3.85
$ws.Range("AX11").Value = 11.25
$ws.Range("AY11").Value = 22
$ws.Range("AZ11").Value = 45
$ws.Range("BA11").Value = 90
$ws.Range("BB11").Value = 350
$ws.Range("G19").Value = 2.55
$ws.Range("I19").Value = 3.1
$ws.Range("AG19").Value = 7.5
$ws.Range("AJ19").Value = 34
$ws.Range("AN19").Value = 4.33
$ws.Range("G29").Value = 2.55
$ws.Range("H29").Value = 3.5
$ws.Range("I29").Value = 2.5
$ws.Range("J29").Value = 3.1
$ws.Range("K29").Value = 2.18
$ws.Range("L29").Value = 3.05
$ws.Range("O29").Value = 1.27
$ws.Range("R29").Value = 1.93
$ws.Range("S29").Value = 1.38
$ws.Range("T29").Value = 2.82
$ws.Range("W29").Value = 9.25
$ws.Range("X29").Value = 13
$ws.Range("Y29").Value = 9.75
$ws.Range("Z29").Value = 27
$ws.Range("AA29").Value = 20
$ws.Range("AB29").Value = 28
$ws.Range("AD29").Value = 6.7
$ws.Range("AG29").Value = 9.25
$ws.Range("AH29").Value = 13
$ws.Range("AI29").Value = 9.5
$ws.Range("AJ29").Value = 26
$ws.Range("AK29").Value = 19.5
$ws.Range("AL29").Value = 27
$ws.Range("AN29").Value = 4.55
$ws.Range("AO29").Value = 13.5
$ws.Range("AP29").Value = 21
$ws.Range("AQ29").Value = 55
$ws.Range("AR29").Value = 90
$ws.Range("AT29").Value = 2.82
$ws.Range("AW29").Value = 4.5
$ws.Range("AX29").Value = 13
$ws.Range("AY29").Value = 20
$ws.Range("AZ29").Value = 55
$ws.Range("AH35").Value = 17
$ws.Range("AT35").Value = 2.27
$ws.Range("AW35").Value = 5.1
$ws.Range("BB35").Value = 450
$ws.Range("I36").Value = 2.22
$ws.Range("G42").Value = 2.07
$ws.Range("H42").Value = 3.25
$ws.Range("I42").Value = 3.25
$ws.Range("J42").Value = 2.7
$ws.Range("K42").Value = 2.1
$ws.Range("L42").Value = 3.85
$ws.Range("N42").Value = 7.2
$ws.Range("R42").Value = 1.82
$ws.Range("T42").Value = 2.67
$ws.Range("U42").Value = 1.75
$ws.Range("V42").Value = 1.98
$ws.Range("W42").Value = 7.7
$ws.Range("X42").Value = 10.25
$ws.Range("Y42").Value = 8.5
$ws.Range("Z42").Value = 19.5
$ws.Range("AA42").Value = 16.5
$ws.Range("AC42").Value = 7.2
$ws.Range("AD42").Value = 6.4
$ws.Range("AE42").Value = 14
$ws.Range("AF42").Value = 65
$ws.Range("AG42").Value = 9.75
$ws.Range("AH42").Value = 17
$ws.Range("AI42").Value = 11.5
$ws.Range("AJ42").Value = 45
$ws.Range("AK42").Value = 30
$ws.Range("AL42").Value = 37
$ws.Range("AN42").Value = 4
$ws.Range("AO42").Value = 10.75
$ws.Range("AP42").Value = 19.5
$ws.Range("AQ42").Value = 40
$ws.Range("AR42").Value = 75
$ws.Range("AT42").Value = 2.67
$ws.Range("AU42").Value = 7.2
$ws.Range("AW42").Value = 5.2
$ws.Range("AX42").Value = 18.5
$ws.Range("AY42").Value = 26
$ws.Range("AZ42").Value = 100
$ws.Range("BA42").Value = 150
